# Update countries & provincias Spain
#
# 1) Swap three pairs of adjacent country-name labels (this is what moves
#    the corresponding entries around in the workbook's shared-string
#    table):
#       - A22/A23:   Portugal / Arabia Saudita  -> Arabia Saudita / Portugal
#       - A193/A194: San Vicente y las Granadinas / Namibia -> Namibia / San Vicente y las Granadinas
#       - A217/A218: San Pedro y Miquelon / Comoras -> Comoras / San Pedro y Miquelon
#
# 2) Update the daily case statistics (columns B:H) for the rows whose
#    numbers changed in this data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country labels (drives the sharedStrings.xml reorder) ---
$ws.Range("A22").Value = "Arabia Saudita"
$ws.Range("A23").Value = "Portugal"

$ws.Range("A193").Value = "Namibia"
$ws.Range("A194").Value = "San Vicente y las Granadinas"

$ws.Range("A217").Value = "Comoras"
$ws.Range("A218").Value = "San Pedro y Miquelon"

# --- Row 22 (now Arabia Saudita) updated stats ---
$ws.Range("B22").Value = 25459
$ws.Range("C22").Value = 1362
$ws.Range("D22").Value = 3765
$ws.Range("E22").Value = 21518
$ws.Range("F22").Value = 117
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 176

# --- Row 23 (now Portugal) updated stats ---
$ws.Range("B23").Value = 25351
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 1671
$ws.Range("E23").Value = 22657
$ws.Range("F23").Value = 150
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 1023

# --- Row 41 (Dinamarca) updated stats ---
$ws.Range("B41").Value = 9407
$ws.Range("C41").Value = 96
$ws.Range("D41").Value = 6889
$ws.Range("E41").Value = 2043
$ws.Range("F41").Value = 60
$ws.Range("G41").Value = 15
$ws.Range("H41").Value = 475

# --- Row 62 (Barein) updated stats ---
$ws.Range("B62").Value = 3273
$ws.Range("C62").Value = 103
$ws.Range("D62").Value = 1567
$ws.Range("E62").Value = 1698

# --- Row 71 (Uzbekistan) updated stats ---
$ws.Range("D71").Value = 1271
$ws.Range("E71").Value = 814

# --- Row 82 (Republica de Macedonia) updated stats ---
$ws.Range("B82").Value = 1506
$ws.Range("C82").Value = 12
$ws.Range("D82").Value = 852
$ws.Range("E82").Value = 572
$ws.Range("F82").Value = 21
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 82
